$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "carlos"
$ws.Range("A4").Value = "juan"
$ws.Range("A5").Value = "luis"
$ws.Range("A6").Value = "marco"
$ws.Range("A2").Value = "daniel"

$ws.Range("A3").Select()
